# Updated day 1 and 2 slides after training
#
# 1) Slide master "date" rectangle: cached field text 14.06.2024 -> 17.06.2024
# 2) Slide 1 "Rectangle 3": drop the " & Malte Fischer" co-presenter suffix
# 3) Slide 4: mark slide as hidden in the slide show (show="0")
# 4) Slide 9 schedule placeholder: drop the trailing "<tab><tab>presenter name"
#    text from the Montag / Dienstag / Mittwoch (17./18./19.06.2024) lines

$p = $ppt.ActivePresentation
$tabChar = [char]9
$ampChar = [char]38

# --- 1) Slide master date field -------------------------------------------------
$master = $p.SlideMaster
$dateShape = $master.Shapes.Item(3)   # "Rectangle 6" holding the cached datetime field
$dateText = $dateShape.TextFrame.TextRange.Text
$dateShape.TextFrame.TextRange.Text = $dateText.Replace("14.06.2024", "17.06.2024")

# --- 2) Slide 1 subtitle: remove " & Malte Fischer" ------------------------------
$slide1 = $p.Slides.Item(1)
$subtitleShape = $slide1.Shapes.Item(2)   # "Rectangle 3"
$subtitleRange = $subtitleShape.TextFrame.TextRange
$subtitleText = $subtitleRange.Text
$ampIdx = $subtitleText.IndexOf($ampChar)
if ($ampIdx -ge 0) {
    $subtitleRange.Text = $subtitleText.Substring(0, $ampIdx).TrimEnd()
}

# --- 3) Slide 4: hide from the slide show ---------------------------------------
$slide4 = $p.Slides.Item(4)
$slide4.SlideShowTransition.Hidden = -1   # msoTrue

# --- 4) Slide 9 schedule list: trim trailing presenter names --------------------
$slide9 = $p.Slides.Item(9)
$scheduleShape = $slide9.Shapes.Item(1)   # "Inhaltsplatzhalter 2"
$scheduleRange = $scheduleShape.TextFrame.TextRange
for ($i = 1; $i -le 3; $i++) {
    $para = $scheduleRange.Paragraphs($i)
    $paraText = $para.Text
    $tabIdx = $paraText.IndexOf($tabChar)
    if ($tabIdx -ge 0) {
        $para.Text = $paraText.Substring(0, $tabIdx)
    }
}
